$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1022.6286
$ws.Range("I15").Value = 1022.6286
$ws.Range("K15").Value = 3067.8858
$ws.Range("M15").Value = -2898.8858
$ws.Range("H19").Value = 1244.5625
$ws.Range("J19").Value = 1322.3572
$ws.Range("L19").Value = 1322.3572
$ws.Range("N19").Value = -1672.3572
$ws.Range("H41").Value = 469.33334
$ws.Range("I41").Value = 83.44444
$ws.Range("K41").Value = 83.44444
$ws.Range("M41").Value = 356.55556
$ws.Range("H58").Value = 537.25
$ws.Range("I58").Value = 537.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1611.75
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -1461.75
$ws.Range("H112").Value = 1709.3334
$ws.Range("J112").Value = 2266.8572
$ws.Range("L112").Value = 6800.571599999999
$ws.Range("N112").Value = -9016.571599999999
$ws.Range("H138").Value = 4099.7656
$ws.Range("I138").Value = 2284.7273
$ws.Range("J138").Value = 6031.9033
$ws.Range("K138").Value = 6854.1819
$ws.Range("L138").Value = 18095.7099
$ws.Range("M138").Value = -1714.1819
$ws.Range("N138").Value = -28375.7099
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 611676.4399999999
$ws.Range("I32").Value = 675770
$ws.Range("K32").Value = 675770
$ws.Range("M32").Value = -675483
$ws.Range("H45").Value = 3118.1177
$ws.Range("I45").Value = 3141.2
$ws.Range("K45").Value = 3141.2
$ws.Range("M45").Value = -2764.2
$ws.Range("H63").Value = 3284.7144
$ws.Range("I63").Value = 2801
$ws.Range("J63").Value = 4494
$ws.Range("K63").Value = 2801
$ws.Range("L63").Value = 4494
$ws.Range("M63").Value = -2115
$ws.Range("N63").Value = -5866
$ws.Range("H66").Value = 3284.7144
$ws.Range("I66").Value = 2801
$ws.Range("J66").Value = 4494
$ws.Range("K66").Value = 14005
$ws.Range("L66").Value = 22470
$ws.Range("M66").Value = -10573
$ws.Range("N66").Value = -29334
$ws.Range("H74").Value = 1981775.5
$ws.Range("I74").Value = 3957249
$ws.Range("J74").Value = 6301.909
$ws.Range("K74").Value = 3957249
$ws.Range("L74").Value = 6301.909
$ws.Range("M74").Value = -3956375
$ws.Range("N74").Value = -8049.909
$ws.Range("H77").Value = 1981775.5
$ws.Range("I77").Value = 3957249
$ws.Range("J77").Value = 6301.909
$ws.Range("K77").Value = 19786245
$ws.Range("L77").Value = 31509.545
$ws.Range("M77").Value = -19781877
$ws.Range("N77").Value = -40245.545
$ws.Range("H97").Value = 1175
$ws.Range("I97").Value = 1175
$ws.Range("K97").Value = 1175
$ws.Range("M97").Value = -679
$ws.Range("H135").Value = 94749.5
$ws.Range("J135").Value = 94749.5
$ws.Range("L135").Value = 94749.5
$ws.Range("N135").Value = -104889.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4625.7144
$ws.Range("I94").Value = 4625.7144
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4625.7144
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -4174.7144
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 738978
$ws.Range("J31").Value = 3874.6316
$ws.Range("L31").Value = 3874.6316
$ws.Range("N31").Value = -4464.631600000001
$ws.Range("H34").Value = 738978
$ws.Range("J34").Value = 3874.6316
$ws.Range("L34").Value = 3874.6316
$ws.Range("N34").Value = -4278.631600000001
$ws.Range("H86").Value = 9364.462
$ws.Range("J86").Value = 10165.523
$ws.Range("L86").Value = 10165.523
$ws.Range("N86").Value = -12411.523
$ws.Range("H89").Value = 9364.462
$ws.Range("J89").Value = 10165.523
$ws.Range("L89").Value = 50827.615
$ws.Range("N89").Value = -62059.615
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 206.16667
$ws.Range("I14").Value = 206.16667
$ws.Range("K14").Value = 618.50001
$ws.Range("M14").Value = -445.50001
$ws.Range("H82").Value = 12483.2
$ws.Range("J82").Value = 13737.5
$ws.Range("L82").Value = 41212.5
$ws.Range("N82").Value = -42024.5
$ws.Range("H85").Value = 12483.2
$ws.Range("J85").Value = 13737.5
$ws.Range("L85").Value = 41212.5
$ws.Range("N85").Value = -44020.5
$ws.Range("H113").Value = 662.7857
$ws.Range("I113").Value = 521.25
$ws.Range("J113").Value = 719.4
$ws.Range("K113").Value = 1563.75
$ws.Range("L113").Value = 2158.2
$ws.Range("M113").Value = 606.25
$ws.Range("N113").Value = -6498.2
$ws.Range("H122").Value = 1344593.5
$ws.Range("I122").Value = 2304542.2
$ws.Range("J122").Value = 665.2
$ws.Range("K122").Value = 20740879.8
$ws.Range("L122").Value = 5986.8
$ws.Range("M122").Value = -20738429.8
$ws.Range("N122").Value = -10886.8
$ws.Range("H132").Value = 1051.4546
$ws.Range("J132").Value = 1398.8334
$ws.Range("L132").Value = 12589.5006
$ws.Range("N132").Value = -17649.5006
$ws.Range("H133").Value = 10792.5
$ws.Range("J133").Value = 15099.2
$ws.Range("L133").Value = 45297.60000000001
$ws.Range("N133").Value = -55417.60000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5997.857
$ws.Range("I80").Value = 2997
$ws.Range("K80").Value = 2997
$ws.Range("M80").Value = -1999
$ws.Range("H83").Value = 5997.857
$ws.Range("I83").Value = 2997
$ws.Range("K83").Value = 14985
$ws.Range("M83").Value = -9993
$ws.Range("H97").Value = 2375
$ws.Range("I97").Value = 1500
$ws.Range("K97").Value = 1500
$ws.Range("M97").Value = -1004
$ws.Range("H102").Value = 3668.5
$ws.Range("I102").Value = 3668.5
$ws.Range("K102").Value = 3668.5
$ws.Range("M102").Value = -2046.5
$ws.Range("H113").Value = 2224.818
$ws.Range("I113").Value = 2046
$ws.Range("K113").Value = 2046
$ws.Range("M113").Value = 124
$ws.Range("H122").Value = 205999.4
$ws.Range("I122").Value = 502499.5
$ws.Range("J122").Value = 8332.666999999999
$ws.Range("K122").Value = 1507498.5
$ws.Range("L122").Value = 24998.001
$ws.Range("M122").Value = -1505048.5
$ws.Range("N122").Value = -29898.001
$ws.Range("H126").Value = 4393.8
$ws.Range("J126").Value = 4666.6665
$ws.Range("L126").Value = 13999.9995
$ws.Range("N126").Value = -18939.9995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3468.6667
$ws.Range("I7").Value = 3694.2856
$ws.Range("K7").Value = 3694.2856
$ws.Range("M7").Value = -3582.2856
$ws.Range("H16").Value = 1819
$ws.Range("I16").Value = 1023.75
$ws.Range("K16").Value = 1023.75
$ws.Range("M16").Value = -853.75
$ws.Range("H40").Value = 5682.875
$ws.Range("I40").Value = 5077.3335
$ws.Range("K40").Value = 5077.3335
$ws.Range("M40").Value = -4941.3335
$ws.Range("H68").Value = 5138.2
$ws.Range("J68").Value = 6000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7498
$ws.Range("H71").Value = 5138.2
$ws.Range("J71").Value = 6000
$ws.Range("L71").Value = 30000
$ws.Range("N71").Value = -37488
$ws.Range("H93").Value = 3035.0667
$ws.Range("I93").Value = 1710.909
$ws.Range("K93").Value = 1710.909
$ws.Range("M93").Value = -462.9090000000001
$ws.Range("H122").Value = 7031.143
$ws.Range("I122").Value = 7414.2
$ws.Range("K122").Value = 22242.6
$ws.Range("M122").Value = -19792.6
$ws.Range("H126").Value = 3468.6667
$ws.Range("I126").Value = 3694.2856
$ws.Range("K126").Value = 11082.8568
$ws.Range("M126").Value = -8612.856800000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 118030
$ws.Range("I4").Value = 118030
$ws.Range("K4").Value = 118030
$ws.Range("M4").Value = -117917
$ws.Range("H81").Value = 88406.914
$ws.Range("I81").Value = 5062.625
$ws.Range("J81").Value = 255095.5
$ws.Range("K81").Value = 10125.25
$ws.Range("L81").Value = 510191
$ws.Range("M81").Value = -9064.25
$ws.Range("N81").Value = -512313
$ws.Range("H84").Value = 88406.914
$ws.Range("I84").Value = 5062.625
$ws.Range("J84").Value = 255095.5
$ws.Range("K84").Value = 50626.25
$ws.Range("L84").Value = 2550955
$ws.Range("M84").Value = -45322.25
$ws.Range("N84").Value = -2561563
$ws.Range("H104").Value = 13000
$ws.Range("J104").Value = 13000
$ws.Range("L104").Value = 13000
$ws.Range("N104").Value = -19988
$ws.Range("H122").Value = 38888.727
$ws.Range("I122").Value = 3611.1428
$ws.Range("K122").Value = 10833.4284
$ws.Range("M122").Value = -8383.428400000001
$ws.Range("H126").Value = 3606.2632
$ws.Range("I126").Value = 2352.8667
$ws.Range("J126").Value = 8306.5
$ws.Range("K126").Value = 7058.6001
$ws.Range("L126").Value = 24919.5
$ws.Range("M126").Value = -4588.6001
$ws.Range("N126").Value = -29859.5
